$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new Price strings in column D (e.g. "597.65", "78.40") look numeric,
# so a plain .Value assignment would let Excel's smart entry silently
# re-interpret them as numbers (dropping trailing zeros, using binary
# float, etc.) instead of keeping the literal text the source data uses.
# Force Text format first so the values land as literal strings, exactly
# like the surrounding (already-text) cells, then drop the temporary
# number format again so no extra cell styling is introduced.
$priceCells = @("D2", "D3", "D5", "D6", "D7", "D9", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D35", "D36", "D37", "D38", "D40", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $priceCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "66.624.04"
$ws.Range("E2").Value = "  -4.31%  "

$ws.Range("D3").Value = "3.447.37"
$ws.Range("E3").Value = "  -4.62%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "597.65"
$ws.Range("E5").Value = "  -4.83%  "

$ws.Range("D6").Value = "146.63"
$ws.Range("E6").Value = "  -7.28%  "

$ws.Range("D7").Value = "3.446.59"
$ws.Range("E7").Value = "  -4.61%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "0.481"
$ws.Range("E9").Value = "  -3.06%  "

$ws.Range("E10").Value = "  -5.26%  "

$ws.Range("D11").Value = "7.39"
$ws.Range("E11").Value = "  +0.03%  "

$ws.Range("D12").Value = "0.421"
$ws.Range("E12").Value = "  -4.36%  "

$ws.Range("D13").Value = "0.0000211"
$ws.Range("E13").Value = "  -7.62%  "

$ws.Range("D14").Value = "4.045.16"
$ws.Range("E14").Value = "  -4.30%  "

$ws.Range("D15").Value = "31.38"
$ws.Range("E15").Value = "  -5.92%  "

$ws.Range("D16").Value = "3.456.29"
$ws.Range("E16").Value = "  -4.43%  "

$ws.Range("D17").Value = "66.743.00"
$ws.Range("E17").Value = "  -4.76%  "

$ws.Range("E18").Value = "  -0.98%  "

$ws.Range("D19").Value = "6.38"
$ws.Range("E19").Value = "  -4.55%  "

$ws.Range("D20").Value = "15.15"
$ws.Range("E20").Value = "  -5.38%  "

$ws.Range("D21").Value = "9.82"
$ws.Range("E21").Value = "  -3.63%  "

$ws.Range("D22").Value = "435.91"
$ws.Range("E22").Value = "  -6.11%  "

$ws.Range("D23").Value = "0.611"
$ws.Range("E23").Value = "  -5.29%  "

$ws.Range("D24").Value = "78.40"
$ws.Range("E24").Value = "  -0.43%  "

$ws.Range("E25").Value = "  +0.01%  "

$ws.Range("D26").Value = "3.599.16"
$ws.Range("E26").Value = "  -4.31%  "

$ws.Range("D27").Value = "5.64"
$ws.Range("E27").Value = "  -5.49%  "

$ws.Range("D28").Value = "0.0000118"
$ws.Range("E28").Value = "  -13.11%  "

$ws.Range("D29").Value = "9.73"
$ws.Range("E29").Value = "  -9.13%  "

$ws.Range("D30").Value = "8.27"
$ws.Range("E30").Value = "  -10.01%  "

$ws.Range("D31").Value = "2.46"
$ws.Range("E31").Value = "  -6.31%  "

$ws.Range("D32").Value = "1.59"
$ws.Range("E32").Value = "  -7.61%  "

$ws.Range("E33").Value = "  +0.27%  "

$ws.Range("E34").Value = "  -8.19%  "

$ws.Range("D35").Value = "25.22"
$ws.Range("E35").Value = "  -5.06%  "

$ws.Range("D36").Value = "6.06"
$ws.Range("E36").Value = "  -7.66%  "

$ws.Range("D37").Value = "3.449.81"
$ws.Range("E37").Value = "  -4.51%  "

$ws.Range("D38").Value = "1.79"
$ws.Range("E38").Value = "  -8.64%  "

$ws.Range("E39").Value = "  -0.02%  "

$ws.Range("D40").Value = "7.83"
$ws.Range("E40").Value = "  -7.47%  "

$ws.Range("E41").Value = "  -0.11%  "

$ws.Range("D42").Value = "2.20"
$ws.Range("E42").Value = "  -8.76%  "

$ws.Range("D43").Value = "173.88"
$ws.Range("E43").Value = "  -3.37%  "

$ws.Range("D44").Value = "0.0880"
$ws.Range("E44").Value = "  -4.92%  "

$ws.Range("D45").Value = "5.33"
$ws.Range("E45").Value = "  -6.15%  "

$ws.Range("D46").Value = "0.881"
$ws.Range("E46").Value = "  -3.67%  "

$ws.Range("D47").Value = "29.67"
$ws.Range("E47").Value = "  -8.73%  "

$ws.Range("D48").Value = "46.11"
$ws.Range("E48").Value = "  +0.17%  "

$ws.Range("D49").Value = "1.23"
$ws.Range("E49").Value = "  -10.08%  "

$ws.Range("D50").Value = "7.45"
$ws.Range("E50").Value = "  -4.74%  "

$ws.Range("D51").Value = "2.40"
$ws.Range("E51").Value = "  -12.87%  "

# Restore default (General) styling on the forced-text cells -- only the
# cell content changed, not its format.
foreach ($addr in $priceCells) { $ws.Range($addr).ClearFormats() }
